$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.726.15'
$ws.Range("D3").Value = '1.633.53'
$ws.Range("E3").Value = '  -0.43%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.29'
$ws.Range("E5").Value = '  -0.29%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.502'
$ws.Range("E6").Value = '  -0.77%  '
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.54'
$ws.Range("E10").Value = '  -4.49%  '
$ws.Range("E11").Value = '  +0.99%  '
$ws.Range("E12").Value = '  -0.87%  '
$ws.Range("D13").Value = '1.857.42'
$ws.Range("E13").Value = '  -0.50%  '
$ws.Range("D14").Value = '1.627.82'
$ws.Range("E14").Value = '  -1.12%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.555'
$ws.Range("E15").Value = '  -1.37%  '
$ws.Range("D16").Value = '0.0₃0766'
$ws.Range("E16").Value = '  -0.13%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '62.64'
$ws.Range("E17").Value = '  -1.24%  '
$ws.Range("D18").Value = '25.740.44'
$ws.Range("E18").Value = '  -0.39%  '
$ws.Range("E19").Value = '  +0.07%  '
$ws.Range("E20").Value = '  +0.87%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '193.16'
$ws.Range("E21").Value = '  +0.15%  '
$ws.Range("E22").Value = '  -0.13%  '
$ws.Range("E23").Value = '  +1.88%  '
$ws.Range("E24").Value = '  +0.05%  '
$ws.Range("E25").Value = '  +2.10%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '140.17'
$ws.Range("E26").Value = '  +0.49%  '
$ws.Range("E27").Value = '  -1.89%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.88'
$ws.Range("E28").Value = '  +0.66%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.48'
$ws.Range("E29").Value = '  -0.83%  '
$ws.Range("E30").Value = '  -0.21%  '
$ws.Range("E31").Value = '  -0.89%  '
$ws.Range("E32").Value = '  +0.92%  '
$ws.Range("E33").Value = '  -0.28%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.59'
$ws.Range("E34").Value = '  +0.80%  '
$ws.Range("E35").Value = '  +0.18%  '
$ws.Range("E36").Value = '  -0.71%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.547'
$ws.Range("E37").Value = '  -1.56%  '
$ws.Range("D38").Value = '1.118.39'
$ws.Range("E38").Value = '  -1.39%  '
$ws.Range("E39").Value = '  -2.02%  '
$ws.Range("E40").Value = '  -0.99%  '
$ws.Range("E41").Value = '  +0.64%  '
$ws.Range("E42").Value = '  +1.03%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '99.66'
$ws.Range("E43").Value = '  +0.59%  '
$ws.Range("E44").Value = '  -0.17%  '
$ws.Range("D45").Value = '1.766.97'
$ws.Range("E45").Value = '  -0.56%  '
$ws.Range("E46").Value = '  -0.08%  '
$ws.Range("E47").Value = '  -1.15%  '
$ws.Range("E48").Value = '  -2.25%  '
$ws.Range("E49").Value = '  -0.43%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.61'
$ws.Range("E50").Value = '  -2.23%  '
$ws.Range("E51").Value = '  +2.81%  '
